$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the start-time / end-time header cells entirely (F1:G1)
$ws.Range("F1:G1").Clear()

# Clear the recorded start/end time values but keep the existing time formatting
$ws.Range("F3:G3").Value = ""
$ws.Range("F5:G5").Value = ""
$ws.Range("F7:G14").Value = ""

# Widen the Description column (E)
$ws.Columns.Item(5).ColumnWidth = 105.67

# Add the new progress entry in row 16
$ws.Range("A16").Value = 4

$ws.Range("B16").Value = 45590
$ws.Range("B16").NumberFormat = "mm-dd-yy"

$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1

$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Lambda Functions"

$ws.Range("E16").Value = "AWS Lambda versions, Aliases, CodeDeploy, Function URL & Security, CodeGuru Profiling, Best Practices"

$ws.Range("F14").Copy()
$ws.Range("F16").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Move the selection to the newly added description cell
$ws.Range("E16").Select()
